# fix(gui) step 1 and 2
# PITONES.xlsx price list: bump the sheet date by one day and refresh the
# "CON TOPE" (step 1, rows 34-37) and "SIN TOPE" (step 2, rows 41-44) prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date stamp in A1: 45308 (2024-01-17) -> 45309 (2024-01-18)
$ws.Range("A1").Value = 45309

# Step 1 - "CON TOPE" price column (D34:D37)
$ws.Range("D34").Value = 1996.418
$ws.Range("D35").Value = 2661.888
$ws.Range("D36").Value = 3855.82
$ws.Range("D37").Value = 3279.202

# Step 2 - "SIN TOPE" price column (D41:D44)
$ws.Range("D41").Value = 2564.02
$ws.Range("D42").Value = 3387.579
$ws.Range("D43").Value = 4266.839
$ws.Range("D44").Value = 3729.338
